# Populate the per-pitch "Pitch / Choice / Result" data (columns F/G/H) that
# feeds the new strikezone visual for hitters, fill in the Exit Velo /
# Launch Angle readings that are now known, fix the at-bat Result labels
# that were placeholders ("Undefined"), and reorder the Pitch Mix lists.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- At-bat block starting row 9 (Inning 3) ----
$ws.Range("F10").Value = "CB"
$ws.Range("G10").Value = "Take"
$ws.Range("H10").Value = "Ball"
$ws.Range("F11").Value = "FB"
$ws.Range("G11").Value = "Take"
$ws.Range("H11").Value = "Strike"
$ws.Range("F12").Value = "FB"
$ws.Range("G12").Value = "Swing"
$ws.Range("H12").Value = "Foul"
$ws.Range("M12").Value = $null
$ws.Range("F13").Value = "CB"
$ws.Range("G13").Value = "Take"
$ws.Range("H13").Value = "Ball"
$ws.Range("F14").Value = "CB"
$ws.Range("G14").Value = "Take"
$ws.Range("H14").Value = "Strike"
$ws.Range("M15").Value = "Strikeout"
$ws.Range("J17").Value = "CH,CB,FB"

# ---- At-bat block starting row 18 (Inning 4) ----
$ws.Range("F19").Value = "CB"
$ws.Range("G19").Value = "Take"
$ws.Range("H19").Value = "Ball"
$ws.Range("M19").Value = "92.3 MPH"
$ws.Range("F20").Value = "FB"
$ws.Range("G20").Value = "Take"
$ws.Range("H20").Value = "Ball"
$ws.Range("F21").Value = "FB"
$ws.Range("G21").Value = "Take"
$ws.Range("H21").Value = "Strike"
$ws.Range("M21").Value = "-3.48°"
$ws.Range("F22").Value = "FB"
$ws.Range("G22").Value = "Swing"
$ws.Range("H22").Value = "In Play"
$ws.Range("J26").Value = "CH,CB,FB"

# ---- At-bat block starting row 27 ----
$ws.Range("F28").Value = "CB"
$ws.Range("G28").Value = "Swing"
$ws.Range("H28").Value = "Strike"
$ws.Range("F29").Value = "SL"
$ws.Range("G29").Value = "Take"
$ws.Range("H29").Value = "Strike"
$ws.Range("F30").Value = "CH"
$ws.Range("G30").Value = "Take"
$ws.Range("H30").Value = "Strike"
$ws.Range("M30").Value = $null
$ws.Range("M33").Value = "Strikeout"
$ws.Range("J35").Value = "CH,CB,FB,SL"

# ---- At-bat block starting row 36 ----
$ws.Range("F37").Value = "CB"
$ws.Range("G37").Value = "Take"
$ws.Range("H37").Value = "Ball"
$ws.Range("M37").Value = "62.56 MPH"
$ws.Range("F38").Value = "FB"
$ws.Range("G38").Value = "Take"
$ws.Range("H38").Value = "Ball"
$ws.Range("F39").Value = "CH"
$ws.Range("G39").Value = "Take"
$ws.Range("H39").Value = "Strike"
$ws.Range("M39").Value = "41.13°"
$ws.Range("F40").Value = "FB"
$ws.Range("G40").Value = "Take"
$ws.Range("H40").Value = "Ball"
$ws.Range("F41").Value = "CH"
$ws.Range("G41").Value = "Swing"
$ws.Range("H41").Value = "In Play"
$ws.Range("J44").Value = "CH,CB,FB,SL"

# ---- At-bat block starting row 45 ----
$ws.Range("F46").Value = "SL"
$ws.Range("G46").Value = "Take"
$ws.Range("H46").Value = "Strike"
$ws.Range("F47").Value = "CH"
$ws.Range("G47").Value = "Swing"
$ws.Range("H47").Value = "Strike"
$ws.Range("F48").Value = "SL"
$ws.Range("G48").Value = "Take"
$ws.Range("H48").Value = "Ball"
$ws.Range("M48").Value = $null
$ws.Range("F49").Value = "SL"
$ws.Range("G49").Value = "Take"
$ws.Range("H49").Value = "Ball"
$ws.Range("F50").Value = "FB"
$ws.Range("G50").Value = "Take"
$ws.Range("H50").Value = "Ball"
$ws.Range("F51").Value = "SL"
$ws.Range("G51").Value = "Take"
$ws.Range("H51").Value = "Ball"
$ws.Range("M51").Value = "Walk"
$ws.Range("J53").Value = "CH,FB,SL"

# ---- At-bat block starting row 60 (Inning 9) ----
$ws.Range("F61").Value = "CH"
$ws.Range("G61").Value = "Swing"
$ws.Range("H61").Value = "In Play"
$ws.Range("M61").Value = "77.19 MPH"
$ws.Range("M63").Value = "-16.33°"
$ws.Range("J68").Value = "CH,FB,SL"
